$d = $word.ActiveDocument

function Get-ParaIndexContaining($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1. Residence changed from "Fort Worth, Texas" to "Colorado Springs"
#    (three runs collapse into a single plain run)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Fort Worth, Texas", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Colorado Springs", 2)

# ---------------------------------------------------------------------------
# 2. Append new "years of experience" entries to the Programming Languages
#    line, mimicking Word's autocorrect/spell-check run-splitting around
#    the newly typed, not-in-dictionary words.
# ---------------------------------------------------------------------------
$skillsIdx = Get-ParaIndexContaining $d "Programming Languages:"
$skillsPara = $d.Paragraphs.Item($skillsIdx)
$skillsEnd = $skillsPara.Range.End - 1
$skillsTarget = $d.Range($skillsEnd - 1, $skillsEnd)   # the trailing ")" run

$skillsBody = `
  '<w:p>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Pytorch</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> (1 </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>yr</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>OpenAI</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> (1 </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>yr</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '</w:p>'
$skillsTarget.InsertXML($pkgHeader + $skillsBody + $pkgFooter)

# ---------------------------------------------------------------------------
# 3. Extend the certifications paragraph with a second course entry. The
#    trailing run that used to read ", Artificial Intelligence Masterclass"
#    is split into three runs (all sharing the original run formatting).
# ---------------------------------------------------------------------------
$certIdx = Get-ParaIndexContaining $d "Artificial Intelligence Masterclass"
$certPara = $d.Paragraphs.Item($certIdx)
$certEnd = $certPara.Range.End - 1
$oldCertText = ", Artificial Intelligence Masterclass"
$certTarget = $d.Range($certEnd - $oldCertText.Length, $certEnd)

$rpr = '<w:rPr><w:rStyle w:val="normaltextrun"/><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr>'
$tm = [char]0x2122
$certBody = `
  '<w:p>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve">, Artificial Intelligence </w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t>Masterclass, Artificial</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> Intelligence A-Z' + $tm + ' 2023: Build an AI with ChatGPT4</w:t></w:r>' + `
  '</w:p>'
$certTarget.InsertXML($pkgHeader + $certBody + $pkgFooter)
